$d = $word.ActiveDocument

$p6 = $d.Paragraphs.Item(6)
$p13 = $d.Paragraphs.Item(13)
$r = $d.Range($p6.Range.Start, $p13.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0" w:before="29" w:line="240" w:lineRule="auto"/><w:ind w:left="0" w:right="-230" w:firstLine="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:b w:val="1"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">B.  Technical Requirements (Nick)</w:t></w:r><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>
<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0" w:before="29" w:line="240" w:lineRule="auto"/><w:ind w:left="1440" w:right="-230" w:hanging="360"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">       While this portion has been heavily covered in section D, the main functional requirements are that 1. The software allows the manager to see the current capacity of the store as determined by local gov’t/business regulations, change the stores allowed capacity in the event that gov’t or business instructions change, and if the store is at capacity, see the estimated wait time for those who would like to enter the store. 2. For the in-store shopper, the software displays the current percent capacity of the store, and if full, the store displays the amount of time they would have to wait to shop in person, making recommendations based on the length of the line to deter Covid/viral spread. C. For the ordering/curbside pick up</w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0" w:before="0" w:beforeAutospacing="0" w:line="240" w:lineRule="auto"/><w:ind w:left="1440" w:right="-230" w:hanging="360"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">        The app should be easy to use for the shopper, either by using an API to allow access via a web browser or downloading the app for their phone. </w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:spacing w:after="0" w:afterAutospacing="0" w:before="0" w:beforeAutospacing="0" w:line="240" w:lineRule="auto"/><w:ind w:left="1440" w:right="-230" w:hanging="360"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Regarding the user interface, this has been discussed previously in prior sections, but it should have a shopping interface implemented either in the browser and communicating with the store’s server via an API or the app itself. In the case of the manager, their portion of the application will be solely implemented via the application and the user interface will be designed for the phone. </w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:spacing w:after="0" w:before="0" w:beforeAutospacing="0" w:line="240" w:lineRule="auto"/><w:ind w:left="1440" w:right="0" w:hanging="360"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:u w:val="none"/></w:rPr></w:pPr><w:sdt><w:sdtPr><w:tag w:val="goog_rdk_0"/></w:sdtPr><w:sdtContent><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:cs="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Once again, depending on whether or not the user is a shopper or manager, the task flow for the manager is as follows Validate Credentials → see current status of store capacity and line → offer the ability to change capacity.</w:t></w:r></w:sdtContent></w:sdt></w:p>
<w:p><w:pPr><w:spacing w:after="0" w:before="29" w:line="240" w:lineRule="auto"/><w:ind w:left="1440" w:right="0" w:firstLine="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:sdt><w:sdtPr><w:tag w:val="goog_rdk_1"/></w:sdtPr><w:sdtContent><w:r><w:rPr><w:rFonts w:ascii="Arial Unicode MS" w:cs="Arial Unicode MS" w:eastAsia="Arial Unicode MS" w:hAnsi="Arial Unicode MS"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">For the shopper, the task flow is as follows sign in→ select whether or not they are ordering online or shopping in store→ in the first case, if they are shopping online, creating a cart, adding items to the cart based on availability → pay for goods → schedule a pick up time → receive a digital receipt. For the second case if the shopper is going in store user credentials → see current capacity and line → receive a recommendation whether they should come to the store now or suggest better times to come. </w:t></w:r></w:sdtContent></w:sdt></w:p>
<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:spacing w:after="0" w:before="81" w:line="329" w:lineRule="auto"/><w:ind w:left="1440" w:right="0" w:hanging="360"/><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">The sole inputs on the manager side are user credentials and allowing the manipulation of a global var that determines story capacity. For output, the store manager receives the current capacity and line.</w:t></w:r></w:p>
<w:p><w:pPr><w:spacing w:after="0" w:before="81" w:line="329" w:lineRule="auto"/><w:ind w:left="1440" w:right="0" w:firstLine="0"/><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">For the shopper, depending on whether in person shopping or curbside pickup, they will require different inputs based on either case. In the case of in person, the only input is user credentials, and the output is the current capacity and a recommendation whether to come to the store or come back later. For the ordering shopper via either a webpage or the app, they will provide user credentials, strings to search for products, ints to select how much of an item they want to buy, credit card information to pay for the goods and services, and a string/selection from a drop-down menu for pick up time. They will receive a receipt, booleans on whether or not the products are available, and a pick up time.</w:t></w:r></w:p>
<w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="7"/></w:numPr><w:spacing w:after="0" w:before="81" w:line="329" w:lineRule="auto"/><w:ind w:left="1440" w:right="2326" w:hanging="360"/><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">There are three cases where API’s will need to be employed in order to achieve full functionality. In the first case, there will need to be a way to access the app in a browser for individuals who may not have a smart-phone. In case two, there will need to be an API to interface with the stores inventory in order to determine whether or not a. the item is in stock and b. if in stock, how much of the item(s) are available, and c. prices for the items. Finally, in order to complete the transaction, a payment API will have to exist to process credit card transactions to pay for the goods. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)

$d.Content.Find.Execute("C.  Acceptance Criteria/Interaction Scenarios (Nick)", $true, $false, $false, $false, $false, $true, 1, $false, "C.  Acceptance Criteria/Interaction Scenarios (Kishan)", 2)
